$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Target cluster for rows 2 and 3 from "ECs" to "Resolving-Mac"
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("D3").Value = "Resolving-Mac"

# Row 2 numeric updates
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.862146333333333
$ws.Range("H2").Value = 8.586439
$ws.Range("I2").Value = 0.9778268096017091
$ws.Range("J2").Value = 0.9778268096017091
$ws.Range("M2").Value = 0.07215833333333334
$ws.Range("N2").Value = 0.216475
$ws.Range("Q2").Value = 0.2065277091694445
$ws.Range("R2").Value = 1.858749382525
$ws.Range("S2").Value = 0.9778268096017091
$ws.Range("T2").Value = 0.9778268096017091

# Row 3 numeric updates
$ws.Range("I3").Value = 0.02217319039829088
$ws.Range("J3").Value = 0.02217319039829088
$ws.Range("M3").Value = 0.07215833333333334
$ws.Range("N3").Value = 0.216475
$ws.Range("Q3").Value = 0.00468322015
$ws.Range("R3").Value = 0.04214898135
$ws.Range("S3").Value = 0.02217319039829088
$ws.Range("T3").Value = 0.02217319039829088
